# Add a new BOM row for the "1X40 Pin Header" / "Theensy / OLED connector" part.
# This inserts a new row at row 45 (pushing the existing rows 45-48 down to
# 46-49) and fills in the Qty/Value/Device columns for the new part.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 45 ("2X5 Pin Header" / POWER).
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45.
$ws.Range("A45").Value = 1
$ws.Range("B45").Value = "1X40 Pin Header"
$ws.Range("C45").Value = "Theensy / OLED connector"

# Match the author's final selection/scroll position.
$ws.Range("D45").Select()
$excel.ActiveWindow.ScrollRow = 27
